$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'57.656.12"
$ws.Range("E2").Value = "  -1.74%  "

# Row 3
$ws.Range("D3").Value = "'3.095.66"
$ws.Range("E3").Value = "  -2.56%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").Value = "'527.90"
$ws.Range("E5").Value = "  -1.52%  "

# Row 6
$ws.Range("D6").Value = "'137.57"
$ws.Range("E6").Value = "  -3.85%  "

# Row 7
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("D8").Value = "'3.098.18"
$ws.Range("E8").Value = "  -2.50%  "

# Row 9
$ws.Range("D9").Value = "'0.468"
$ws.Range("E9").Value = "  +4.18%  "

# Row 10
$ws.Range("D10").Value = "'7.26"
$ws.Range("E10").Value = "  +0.34%  "

# Row 11
$ws.Range("D11").Value = "'0.107"
$ws.Range("E11").Value = "  -3.54%  "

# Row 12
$ws.Range("D12").Value = "'0.406"
$ws.Range("E12").Value = "  +1.45%  "

# Row 13
$ws.Range("E13").Value = "  +1.88%  "

# Row 14
$ws.Range("D14").Value = "'3.614.51"
$ws.Range("E14").Value = "  -3.07%  "

# Row 15
$ws.Range("D15").Value = "'25.43"
$ws.Range("E15").Value = "  -2.27%  "

# Row 16
$ws.Range("D16").Value = "'0.0000162"
$ws.Range("E16").Value = "  -3.43%  "

# Row 17
$ws.Range("D17").Value = "'57.671.14"
$ws.Range("E17").Value = "  -1.85%  "

# Row 18
$ws.Range("D18").Value = "'3.082.66"
$ws.Range("E18").Value = "  -2.99%  "

# Row 19
$ws.Range("D19").Value = "'5.93"
$ws.Range("E19").Value = "  -4.36%  "

# Row 20
$ws.Range("D20").Value = "'12.55"
$ws.Range("E20").Value = "  -3.24%  "

# Row 21
$ws.Range("D21").Value = "'7.89"
$ws.Range("E21").Value = "  -2.54%  "

# Row 22
$ws.Range("D22").Value = "'350.67"
$ws.Range("E22").Value = "  -1.88%  "

# Row 23
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.37%  "

# Row 24
$ws.Range("D24").Value = "'68.48"
$ws.Range("E24").Value = "  -0.16%  "

# Row 25
$ws.Range("D25").Value = "'0.501"
$ws.Range("E25").Value = "  -3.17%  "

# Row 26
$ws.Range("D26").Value = "'0.167"
$ws.Range("E26").Value = "  -2.11%  "

# Row 27
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.27%  "

# Row 28
$ws.Range("D28").Value = "'0.0₃0866"
$ws.Range("E28").Value = "  -9.46%  "

# Row 29
$ws.Range("E29").Value = "  -0.07%  "

# Row 30
$ws.Range("D30").Value = "'7.16"
$ws.Range("E30").Value = "  -5.63%  "

# Row 31
$ws.Range("D31").Value = "'1.86"
$ws.Range("E31").Value = "  -2.89%  "

# Row 32
$ws.Range("D32").Value = "'6.01"
$ws.Range("E32").Value = "  -8.52%  "

# Row 33
$ws.Range("D33").Value = "'21.18"
$ws.Range("E33").Value = "  -1.22%  "

# Row 34
$ws.Range("D34").Value = "'4.87"
$ws.Range("E34").Value = "  -1.04%  "

# Row 35
$ws.Range("D35").Value = "'159.43"
$ws.Range("E35").Value = "  +0.96%  "

# Row 36
$ws.Range("D36").Value = "'1.14"
$ws.Range("E36").Value = "  -6.81%  "

# Row 37
$ws.Range("D37").Value = "'6.03"
$ws.Range("E37").Value = "  -3.95%  "

# Row 38
$ws.Range("D38").Value = "'25.77"
$ws.Range("E38").Value = "  -3.12%  "

# Row 39
$ws.Range("D39").Value = "'1.25"
$ws.Range("E39").Value = "  -5.44%  "

# Row 40
$ws.Range("D40").Value = "'0.0665"
$ws.Range("E40").Value = "  -1.94%  "

# Row 41
$ws.Range("D41").Value = "'1.60"
$ws.Range("E41").Value = "  -4.18%  "

# Row 42
$ws.Range("D42").Value = "'4.02"
$ws.Range("E42").Value = "  -1.09%  "

# Row 43
$ws.Range("D43").Value = "'0.693"
$ws.Range("E43").Value = "  -2.36%  "

# Row 44
$ws.Range("D44").Value = "'2.396.44"
$ws.Range("E44").Value = "  +1.90%  "

# Row 45
$ws.Range("D45").Value = "'36.91"
$ws.Range("E45").Value = "  -0.41%  "

# Row 46
$ws.Range("E46").Value = "  +0.12%  "

# Row 47
$ws.Range("D47").Value = "'3.139.36"
$ws.Range("E47").Value = "  -2.48%  "

# Row 48
$ws.Range("D48").Value = "'0.0265"
$ws.Range("E48").Value = "  -3.39%  "

# Row 49
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").Value = "'0.956"
$ws.Range("E49").Value = "  -6.64%  "

# Row 50
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "'6.03"
$ws.Range("E50").Value = "  -1.28%  "

# Row 51
$ws.Range("D51").Value = "'19.55"
$ws.Range("E51").Value = "  -5.84%  "
